{"js": "const replacements = [\n  [\"544\u00d74=2176\", \"548\u00d77=3836\"],\n  [\"476\u00d74=1904\", \"904\u00d79=8136\"],\n  [\"766\u00d72=1532\", \"155\u00d75=775\"],\n  [\"559\u00d73=1677\", \"417\u00d75=2085\"],\n  [\"568\u00d77=3976\", \"379\u00d75=1895\"],\n  [\"855\u00d75=4275\", \"952\u00d73=2856\"],\n  [\"817\u00d72=1634\", \"836\u00d77=5852\"],\n  [\"835\u00d73=2505\", \"917\u00d78=7336\"],\n  [\"263\u00d72=526\", \"728\u00d73=2184\"],\n  [\"754\u00d76=4524\", \"135\u00d75=675\"],\n  [\"687\u00d74=2748\", \"396\u00d72=792\"],\n  [\"119\u00d72=238\", \"963\u00d73=2889\"],\n  [\"850\u00d74=3400\", \"134\u00d79=1206\"],\n  [\"413\u00d73=1239\", \"614\u00d76=3684\"],\n  [\"508\u00d76=3048\", \"224\u00d73=672\"],\n  [\"897\u00d72=1794\", \"675\u00d76=4050\"],\n  [\"614\u00d75=3070\", \"844\u00d77=5908\"],\n  [\"288\u00d72=576\", \"212\u00d72=424\"],\n  [\"565\u00d72=1130\", \"354\u00d77=2478\"],\n  [\"929\u00d78=7432\", \"794\u00d75=3970\"],\n  [\"471\u00d78=3768\", \"621\u00d74=2484\"],\n  [\"954\u00d73=2862\", \"543\u00d79=4887\"],\n  [\"652\u00d77=4564\", \"272\u00d79=2448\"],\n  [\"812\u00d73=2436\", \"678\u00d75=3390\"],\n  [\"636\u00d77=4452\", \"782\u00d75=3910\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"544\u00d74=2176\"; New = \"548\u00d77=3836\" },\n    @{ Old = \"476\u00d74=1904\"; New = \"904\u00d79=8136\" },\n    @{ Old = \"766\u00d72=1532\"; New = \"155\u00d75=775\" },\n    @{ Old = \"559\u00d73=1677\"; New = \"417\u00d75=2085\" },\n    @{ Old = \"568\u00d77=3976\"; New = \"379\u00d75=1895\" },\n    @{ Old = \"855\u00d75=4275\"; New = \"952\u00d73=2856\" },\n    @{ Old = \"817\u00d72=1634\"; New = \"836\u00d77=5852\" },\n    @{ Old = \"835\u00d73=2505\"; New = \"917\u00d78=7336\" },\n    @{ Old = \"263\u00d72=526\"; New = \"728\u00d73=2184\" },\n    @{ Old = \"754\u00d76=4524\"; New = \"135\u00d75=675\" },\n    @{ Old = \"687\u00d74=2748\"; New = \"396\u00d72=792\" },\n    @{ Old = \"119\u00d72=238\"; New = \"963\u00d73=2889\" },\n    @{ Old = \"850\u00d74=3400\"; New = \"134\u00d79=1206\" },\n    @{ Old = \"413\u00d73=1239\"; New = \"614\u00d76=3684\" },\n    @{ Old = \"508\u00d76=3048\"; New = \"224\u00d73=672\" },\n    @{ Old = \"897\u00d72=1794\"; New = \"675\u00d76=4050\" },\n    @{ Old = \"614\u00d75=3070\"; New = \"844\u00d77=5908\" },\n    @{ Old = \"288\u00d72=576\"; New = \"212\u00d72=424\" },\n    @{ Old = \"565\u00d72=1130\"; New = \"354\u00d77=2478\" },\n    @{ Old = \"929\u00d78=7432\"; New = \"794\u00d75=3970\" },\n    @{ Old = \"471\u00d78=3768\"; New = \"621\u00d74=2484\" },\n    @{ Old = \"954\u00d73=2862\"; New = \"543\u00d79=4887\" },\n    @{ Old = \"652\u00d77=4564\"; New = \"272\u00d79=2448\" },\n    @{ Old = \"812\u00d73=2436\"; New = \"678\u00d75=3390\" },\n    @{ Old = \"636\u00d77=4452\"; New = \"782\u00d75=3910\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.Text = $r.New\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Execute($r.Old, $false, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)\n}\n"}
